$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Line spacing: every paragraph's "atLeast" line spacing goes from 14pt
#    (w:line="280") to 15pt (w:line="300").  The spacing rule itself (AtLeast)
#    stays the same, only the minimum value increases.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.ParagraphFormat.LineSpacing = 15
}

# ---------------------------------------------------------------------------
# 2) Replace the hard-coded "HIJO(A)." wording with the new templated
#    "${menorhijo}." placeholder so the document can show more than one
#    child.  Done in two Find/Replace passes so the existing
#    gramStart/gramEnd proof-reading marks (which sit between the two
#    original runs) stay anchored exactly where they were.
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("MI MENOR HIJO(A", $true, $false, $false, $false, $false, $true, 1, $false, "`${menorhijo", 2)

$r2 = $d.Content
$r2.Find.Execute(").", $true, $false, $false, $false, $false, $true, 1, $false, "}.", 2)

# ---------------------------------------------------------------------------
# 3) Remove the single stray space that was sitting in the middle of the
#    "====...== ==...====" divider line right after the paragraph above.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$found3 = $r3.Find.Execute("=========================== ==========")
if ($found3) {
    $spaceStart = $r3.Start + 27
    $spaceEnd = $spaceStart + 1
    $spaceRange = $d.Range($spaceStart, $spaceEnd)
    $spaceRange.Delete()
}

# ---------------------------------------------------------------------------
# 4) Underline the "AUTORIZACION DE VIAJE:" heading run.
# ---------------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("AUTORIZACION DE VIAJE:")
$r4.Font.Underline = 1
